$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (pushes existing row 2 data down to row 3)
$ws.Rows.Item(2).Insert()
# Strip the formatting the insert copied down from the header row above
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new book record
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "67890"
$ws.Cells.Item(2, 1).ClearFormats()

$ws.Cells.Item(2, 2).Value = "Charles Darwin"
$ws.Cells.Item(2, 3).Value = "The Descent of Man"
$ws.Cells.Item(2, 4).Value = "and Selection in Relation to Sex"
$ws.Cells.Item(2, 5).Value = "Berlin"
$ws.Cells.Item(2, 6).Value = "John Murray"

$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "1871"
$ws.Cells.Item(2, 7).ClearFormats()

$ws.Cells.Item(2, 8).Value = "Second edition"
$ws.Cells.Item(2, 9).Value = "ALT-123"
